# Version 2.0.1 solucionado error espera de base de datos
# Update patient record fields on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Patient name and record number
$ws.Range("A6").Value = "GUTIERREZ  RODRIGUEZ  ISIDRO  ANTONIO"
$ws.Range("G6").Value = "/201761926"

# Date of birth and age (leading apostrophe forces text, avoiding
# Excel's automatic date/number reinterpretation of these values)
$ws.Range("A9").Value = "'1971-10-09"
$ws.Range("D9").Value = "'46"

# Sex
$ws.Range("G9").Value = "MASCULINO"

# Occupation, nationality, identification document
$ws.Range("C11").Value = "N.T."
$ws.Range("E11").Value = "GUATEMALTECO"
$ws.Range("G11").Value = "'1946786950101"

# Emergency contact info
$ws.Range("A13").Value = "ESPERANZA CASTILLO"
$ws.Range("D13").Value = "ESPOSA"
$ws.Range("E13").Value = "24 calle 3-28 z. 3"
$ws.Range("G13").Value = "'58772234"

# Time of medical assistance
$ws.Range("D14").Value = "Hora: 11:45:49"
